$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Greetings young one..." introduction text in B2 to mention
# that Medialogy is a bachelor and masters degree at Aalborg University.
$ws.Range("B2").Value = " Greetings young one, and welcome to Medialogy,  a bachelor and masters degree at Aalborg University _To begin your quest, place your finger on top of your avatar and drag yourself through the path of Medialogy."

# Move the active selection from E18 to B16.
$ws.Range("B16").Select()
